# Weekly price-sheet update: a new daily record is inserted right after the
# existing row 26, pushing all subsequent records down by one row (so the
# former row 67 becomes row 68). The new record carries the values that the
# diff shows for the new row 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 26; Excel shifts rows 26-67 down to 27-68
# (including formatting / styles), growing the used range to A1:R68.
$ws.Rows.Item(26).Insert()

# Populate the newly-inserted row 26 with the new observation's data.
$ws.Cells.Item(26, 1).Value  = 10
$ws.Cells.Item(26, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value  = "La Araucanía"
$ws.Cells.Item(26, 4).Value  = 44721
$ws.Cells.Item(26, 5).Value  = 9
$ws.Cells.Item(26, 6).Value  = 300000001
$ws.Cells.Item(26, 7).Value  = "Rabanito"
$ws.Cells.Item(26, 8).Value  = "Sin especificar"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 7000
$ws.Cells.Item(26, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(26, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(26, 16).Value = 583
$ws.Cells.Item(26, 17).Value = 12
$ws.Cells.Item(26, 18).Value = "Hortaliza"
